$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. Resize the table grid columns (tblGrid/gridCol widths), in twentieths
#    of a point (dxa). The 6th column already matches and is left untouched.
$newWidths = @(3724, 2306, 1937, 3944, 4342)
for ($i = 0; $i -lt $newWidths.Length; $i++) {
    $t.Columns.Item($i + 1).Width = $newWidths[$i] / 20.0
}

# 2. Row 3 ("Across all sessions..."), Actual Output column: append more text
#    explaining the exception.
$cell = $t.Cell(3, 5)
$p = $cell.Range.Paragraphs.Item(1)
$r = $p.Range
$endRange = $d.Range($r.End - 1, $r.End - 1)
$endRange.InsertAfter(" except when accessing the statistics screen")

# 3. Row 3, Pass/Fail column: change the verdict from Pass to Fail.
#    NB: Table-cell Range objects don't scope .Find correctly in this
#    runtime (it searches from the top of the document), so re-wrap the
#    cell's Start/End into a plain Document.Range, and use wdReplaceOne
#    (not wdReplaceAll) so only the single match inside that span is hit.
$cell = $t.Cell(3, 6)
$cr = $cell.Range
$rng = $d.Range($cr.Start, $cr.End)
$rng.Find.Execute("Pass", $true, $false, $false, $false, $false, $true, 0, $false, "Fail", 1)

# 4. Row 4 ("First user's settings..."), Actual Output column: normalize the
#    run/proofErr structure into a single plain run (text content unchanged).
$cell = $t.Cell(4, 5)
$cr = $cell.Range
$rng = $d.Range($cr.Start, $cr.End)
$quote = [char]0x2019
$targetText = "First user" + $quote + "s settings carries over to the second user" + $quote + "s background settings"
$rng.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 0, $false, $targetText, 1)
